# CF1: Divide adjacent repeat rules. Fix instruments
#
# The "Bass" header labels in row 3 are renamed to "Cantus" (C3/G3).
# (Rows 5-7 in column A keep their existing text; only the shared-string
# table layout shifted around them because the old "No bass change" /
# "Bass change" strings were dropped and replaced by the new ones.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels in row 3 (C3 / G3): rename Bass -> Cantus
$ws.Range("C3").Value = "No cantus change"
$ws.Range("G3").Value = "Cantus change"

# Update the active selection to match the author's final cursor position.
$ws.Range("G4").Select()
